# Scheduled-runner refresh: update cached market/profit figures across
# the per-job Leve profit sheets (ALC, ARM, BSM, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3156.796
$ws.Range("J17").Value = 3516.262
$ws.Range("L17").Value = 10548.786
$ws.Range("N17").Value = -10884.786

$ws.Range("H40").Value = 6262.727
$ws.Range("I40").Value = 4316.5
$ws.Range("J40").Value = 8598.200000000001
$ws.Range("K40").Value = 4316.5
$ws.Range("L40").Value = 8598.200000000001
$ws.Range("M40").Value = -4141.5
$ws.Range("N40").Value = -8948.200000000001

$ws.Range("H74").Value = 7385.1904
$ws.Range("I74").Value = 4577.5
$ws.Range("K74").Value = 4577.5
$ws.Range("M74").Value = -3641.5

$ws.Range("H77").Value = 7385.1904
$ws.Range("I77").Value = 4577.5
$ws.Range("K77").Value = 22887.5
$ws.Range("M77").Value = -18207.5

$ws.Range("H86").Value = 3661.8572
$ws.Range("I86").Value = 3600
$ws.Range("J86").Value = 3672.1667
$ws.Range("K86").Value = 3600
$ws.Range("L86").Value = 3672.1667
$ws.Range("M86").Value = -2477
$ws.Range("N86").Value = -5918.1667

$ws.Range("H89").Value = 3661.8572
$ws.Range("I89").Value = 3600
$ws.Range("J89").Value = 3672.1667
$ws.Range("K89").Value = 18000
$ws.Range("L89").Value = 18360.8335
$ws.Range("M89").Value = -12384
$ws.Range("N89").Value = -29592.8335

$ws.Range("H112").Value = 3157.2856
$ws.Range("J112").Value = 3483.6667
$ws.Range("L112").Value = 10451.0001
$ws.Range("N112").Value = -12667.0001

$ws.Range("H132").Value = 1799.75
$ws.Range("I132").Value = 1486.973
$ws.Range("K132").Value = 4460.919
$ws.Range("M132").Value = -1930.919

$ws.Range("H138").Value = 5030.316
$ws.Range("I138").Value = 3597.818
$ws.Range("J138").Value = 7000
$ws.Range("K138").Value = 10793.454
$ws.Range("L138").Value = 21000
$ws.Range("M138").Value = -5653.454000000002
$ws.Range("N138").Value = -31280

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 32998
$ws.Range("J24").Value = 32998
$ws.Range("L24").Value = 32998
$ws.Range("N24").Value = -33746

$ws.Range("H88").Value = 5357
$ws.Range("J88").Value = 3960.111
$ws.Range("L88").Value = 3960.111
$ws.Range("N88").Value = -4772.111

$ws.Range("H91").Value = 5357
$ws.Range("J91").Value = 3960.111
$ws.Range("L91").Value = 3960.111
$ws.Range("N91").Value = -6768.111

$ws.Range("H100").Value = 32998
$ws.Range("J100").Value = 32998
$ws.Range("L100").Value = 32998
$ws.Range("N100").Value = -35162

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 42000
$ws.Range("J9").Value = 42000
$ws.Range("L9").Value = 42000
$ws.Range("N9").Value = -42336

$ws.Range("H105").Value = 4805.2144
$ws.Range("I105").Value = 2228.087
$ws.Range("J105").Value = 16660
$ws.Range("K105").Value = 2228.087
$ws.Range("L105").Value = 16660
$ws.Range("M105").Value = -481.087
$ws.Range("N105").Value = -20154

$ws.Range("H107").Value = 2666.3333
$ws.Range("I107").Value = 2666.3333
$ws.Range("K107").Value = 2666.3333
$ws.Range("M107").Value = -746.3332999999998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 9913709
$ws.Range("J122").Value = 10991494
$ws.Range("L122").Value = 98923446
$ws.Range("N122").Value = -98928346

$ws.Range("H129").Value = 5955412.5
$ws.Range("I129").Value = 712.1111
$ws.Range("J129").Value = 16673873
$ws.Range("K129").Value = 2136.3333
$ws.Range("L129").Value = 50021619
$ws.Range("M129").Value = 2863.6667
$ws.Range("N129").Value = -50031619

$ws.Range("H137").Value = 70779.07000000001
$ws.Range("I137").Value = 1224.8
$ws.Range("J137").Value = 105556.2
$ws.Range("K137").Value = 3674.4
$ws.Range("L137").Value = 316668.6
$ws.Range("M137").Value = 1425.6
$ws.Range("N137").Value = -326868.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 29500
$ws.Range("J26").Value = 29500
$ws.Range("L26").Value = 29500
$ws.Range("N26").Value = -30060

$ws.Range("H42").Value = 50000
$ws.Range("J42").Value = 50000
$ws.Range("L42").Value = 50000
$ws.Range("N42").Value = -50970

$ws.Range("H50").Value = 29500
$ws.Range("J50").Value = 29500
$ws.Range("L50").Value = 29500
$ws.Range("N50").Value = -30496

$ws.Range("H52").Value = 25011000

$ws.Range("H93").Value = 20547.785
$ws.Range("J93").Value = 21915.7
$ws.Range("L93").Value = 21915.7
$ws.Range("N93").Value = -25659.7

$ws.Range("H97").Value = 1670.4
$ws.Range("I97").Value = 1188.091
$ws.Range("J97").Value = 2996.75
$ws.Range("K97").Value = 1188.091
$ws.Range("L97").Value = 2996.75
$ws.Range("M97").Value = -692.0909999999999
$ws.Range("N97").Value = -3988.75

$ws.Range("H115").Value = 50000
$ws.Range("J115").Value = 50000
$ws.Range("L115").Value = 50000
$ws.Range("N115").Value = -52350

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = $null

$ws.Range("H126").Value = 3749.077
$ws.Range("I126").Value = 1960.5714
$ws.Range("K126").Value = 5881.7142
$ws.Range("M126").Value = -3411.7142

$ws.Range("H132").Value = 288726.38
$ws.Range("I132").Value = 348392.28
$ws.Range("K132").Value = 1045176.84
$ws.Range("M132").Value = -1042646.84

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 24993
$ws.Range("I23").Value = 24993
$ws.Range("K23").Value = 24993
$ws.Range("M23").Value = -24763

$ws.Range("H93").Value = 512.3
$ws.Range("I93").Value = 518.25
$ws.Range("K93").Value = 518.25
$ws.Range("M93").Value = 729.75

$ws.Range("H122").Value = 4490.846
$ws.Range("I122").Value = 2826.7
$ws.Range("J122").Value = 10038
$ws.Range("K122").Value = 8480.099999999999
$ws.Range("L122").Value = 30114
$ws.Range("M122").Value = -6030.099999999999
$ws.Range("N122").Value = -35014

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 14374.875
$ws.Range("J14").Value = 5000
$ws.Range("L14").Value = 5000
$ws.Range("N14").Value = -5336

$ws.Range("H70").Value = 36332.5
$ws.Range("I70").Value = 33999.25
$ws.Range("K70").Value = 33999.25
$ws.Range("M70").Value = -33684.25

$ws.Range("H73").Value = 36332.5
$ws.Range("I73").Value = 33999.25
$ws.Range("K73").Value = 33999.25
$ws.Range("M73").Value = -32907.25

$ws.Range("H93").Value = 52996.668
$ws.Range("J93").Value = 52996.668
$ws.Range("L93").Value = 52996.668
$ws.Range("N93").Value = -57988.668

$ws.Range("H122").Value = 3437.6
$ws.Range("I122").Value = 1708.1666
$ws.Range("J122").Value = 19002.5
$ws.Range("K122").Value = 5124.4998
$ws.Range("L122").Value = 57007.5
$ws.Range("M122").Value = -2674.4998
$ws.Range("N122").Value = -61907.5
